$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: AccountDetails block ---
$ws.Range("A2").Value = "AccountDetails"

$ws.Range("B2").Value = "vnarra@helenoftroy.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:vnarra@helenoftroy.com")

$ws.Range("C2").Value = "avayugundla@helenoftroy.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com")

$ws.Range("D2").Value = "Lotuswave@123"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123")

$ws.Range("E2").Value = "Lotuswave@123"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123")

$ws.Range("F2").Value = "Test"
$ws.Range("G2").Value = "Qa"

# --- Row 3: Product block ---
$ws.Range("A3").Value = "Product"
$ws.Range("X3").Value = "Curl Defining Styling Soufflé (Mini 2 fl oz.)"
$ws.Range("Y3").Value = "'1"
$ws.Range("AA3").Value = 1

# --- Row 4: Configurable Product block ---
$ws.Range("A4").Value = "Configurable Product"
$ws.Range("X4").Value = "Curl Defining Styling Soufflé (Mini 2 fl oz.)"
$ws.Range("Y4").Value = "'1"
$ws.Range("AA4").Value = 1
$ws.Range("AB4").Value = "Mini 2 fl oz."
$ws.Range("AB4").Font.Name = "Courier New"
$ws.Range("AB4").Font.Size = 9
$ws.Range("AB4").Font.Color = 2039583

# --- Row 5: Address block ---
$ws.Range("A5").Value = "Address"
$ws.Range("F5").Value = "Qa"
$ws.Range("G5").Value = "Test"

$ws.Range("K5").Value = "vnarra@helenoftroy.com"
$ws.Hyperlinks.Add($ws.Range("K5"), "mailto:vnarra@helenoftroy.com")

$ws.Range("N5").Value = "844 N Colony Rd"
$ws.Range("O5").Value = "Wallingford"
$ws.Range("P5").Value = "United States"
$ws.Range("Q5").Value = "Connecticut"
$ws.Range("R5").Value = "'06492"
$ws.Range("S5").Value = 9898989898

# --- View state ---
$ws.Range("K11").Select()
$excel.ActiveWindow.ScrollColumn = 10
